# Update "想去人数" (F column) figures for several events that appear
# on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 108   # 南宁·原神x星铁x绝区零同人ONLY3.0: 107 -> 108
$wsExpo.Range("F5").Value = 2881  # 南宁·2024良牙动漫秋季盛典（秋典）: 2868 -> 2881
$wsExpo.Range("F6").Value = 288   # 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini: 287 -> 288
$wsExpo.Range("F7").Value = 393   # 南宁·万圣漫控嘉年华10: 392 -> 393

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 108    # 南宁·原神x星铁x绝区零同人ONLY3.0: 107 -> 108
$wsAll.Range("F5").Value = 2881   # 南宁·2024良牙动漫秋季盛典（秋典）: 2868 -> 2881
$wsAll.Range("F6").Value = 288    # 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini: 287 -> 288
$wsAll.Range("F9").Value = 393    # 南宁·万圣漫控嘉年华10: 392 -> 393
